$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5244542360305786
$ws.Range("B1").Value = 2.37044358253479
$ws.Range("C1").Value = 6.561444282531738
$ws.Range("D1").Value = 1.988759636878967
$ws.Range("E1").Value = 1.632404088973999
